$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update cell text content ---
# Note: "view history a single item" must be written before the updated
# "outgoing/incoming" text so the new entries land in the shared-strings
# table in the same order as the target file.
$ws.Range("D7").Value = "view history a single item"
$ws.Range("D6").Value = "View detailed list of orders that are outgoing/incoming"

# --- Give D7 the same formatting (font) already used by the D10:D15 block (style index 2) ---
$ws.Range("D10").Copy() | Out-Null
$ws.Range("D7").PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats)
$excel.CutCopyMode = 0

# --- Update the active selection shown in the sheet view ---
$ws.Range("C10").Select() | Out-Null
